$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.8656368598586758
$ws.Range("D2").Value = 0.3927566250147447

$ws.Range("C3").Value = 0.9862376393140799
$ws.Range("D3").Value = 0.3309833079280751

$ws.Range("C4").Value = -0.9297158228684781
$ws.Range("D4").Value = 0.3590725484231903

$ws.Range("C5").Value = -1.548317702376349
$ws.Range("D5").Value = 0.1308047655289333
$ws.Range("G5").Value = "No"

$ws.Range("C6").Value = 0.01960575829442326
$ws.Range("D6").Value = 0.9844724769138642

$ws.Range("C7").Value = -1.33266145232974
$ws.Range("D7").Value = 0.1915056314261376

$ws.Range("C8").Value = -1.510536943227026
$ws.Range("D8").Value = 0.1401453774405939
$ws.Range("G8").Value = "No"

$ws.Range("C9").Value = -1.528878351364493
$ws.Range("D9").Value = 0.1355458652578285

$ws.Range("C10").Value = -1.524384082660636
$ws.Range("D10").Value = 0.1366614875016696
$ws.Range("G10").Value = "No"

$ws.Range("C11").Value = -1.312714718197064
$ws.Range("D11").Value = 0.1980713222591872

$wb.Save()
